$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row (row 1)
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Style the header row to match the other sheets: bold font, thin border,
# centered horizontally, aligned to top
$headerRange = $ws.Range("A1:F1")
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Font.Bold = $true

# Data row (row 2)
# Force text storage (not numeric) for the match code, matching the rest of
# the workbook where such codes are stored as text
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4656"

# These columns have no data for this row, but the cells themselves are
# still present (empty) in the source sheet. Touching a formatting property
# that is already at its default value forces the engine to keep the cell
# entry without altering its (default) style.
$ws.Range("B2:E2").Interior.Pattern = -4142

$ws.Range("F2").Value = "NO"
